$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3014822451252854
$ws.Range("D2").Value = 0.0008071150410611239
$ws.Range("E2").Value = 0.08742300873016973
$ws.Range("F2").Value = 6.488479599838342
$ws.Range("G2").Value = 0.002740648245349902
$ws.Range("J2").Value = 0.2869176022684456
$ws.Range("K2").Value = 5.83696212664961
$ws.Range("L2").Value = 0.05403137655040613
$ws.Range("M2").Value = 1.320918520843229
$ws.Range("N2").Value = 3.626605083532354

$ws.Range("C3").Value = 0.3004137547442696
$ws.Range("D3").Value = 0.000721096781036934
$ws.Range("E3").Value = 0.08769605131030289
$ws.Range("F3").Value = 6.477801074685971
$ws.Range("G3").Value = 0.002747390666001467
$ws.Range("J3").Value = 0.2879893418097836
$ws.Range("K3").Value = 5.723817180787705
$ws.Range("L3").Value = 0.05395297221628415
$ws.Range("M3").Value = 1.303538214084689
$ws.Range("N3").Value = 3.649992641181726

$ws.Range("C4").Value = 0.2998917496791051
$ws.Range("D4").Value = 0.0006689668209638455
$ws.Range("E4").Value = 0.08788727671919538
$ws.Range("F4").Value = 6.474164032545417
$ws.Range("G4").Value = 0.002751746420529888
$ws.Range("J4").Value = 0.2887389460800627
$ws.Range("K4").Value = 5.65773650135236
$ws.Range("L4").Value = 0.05390451724019218
$ws.Range("M4").Value = 1.293579438982093
$ws.Range("N4").Value = 3.665271115587295

$ws.Range("C5").Value = 0.299712717706754
$ws.Range("D5").Value = 0.0006478725100969029
$ws.Range("E5").Value = 0.08797114313364673
$ws.Range("F5").Value = 6.473414264369396
$ws.Range("G5").Value = 0.002753575908728324
$ws.Range("J5").Value = 0.2890674365232258
$ws.Range("K5").Value = 5.631658163607767
$ws.Range("L5").Value = 0.05388468814648739
$ws.Range("M5").Value = 1.28970004619292
$ws.Range("N5").Value = 3.671727711005829

$ws.Range("C6").Value = 0.2996850238190234
$ws.Range("D6").Value = 0.0006443781424021466
$ws.Range("E6").Value = 0.08798542823874822
$ws.Range("F6").Value = 6.473333946446928
$ws.Range("G6").Value = 0.002753882990171244
$ws.Range("J6").Value = 0.2891233724933855
$ws.Range("K6").Value = 5.627379144716656
$ws.Range("L6").Value = 0.05388139038726347
$ws.Range("M6").Value = 1.289066670462766
$ws.Range("N6").Value = 3.672813731896888

$ws.Range("C7").Value = 0.2998891988066106
$ws.Range("D7").Value = 0.0006686817633188014
$ws.Range("E7").Value = 0.08788838370181473
$ws.Range("F7").Value = 6.474150957959026
$ws.Range("G7").Value = 0.002751770872572569
$ws.Range("J7").Value = 0.2887432830070047
$ws.Range("K7").Value = 5.657381361148111
$ws.Range("L7").Value = 0.05390425016095257
$ws.Range("M7").Value = 1.293526396267765
$ws.Range("N7").Value = 3.665357258814858

$ws.Range("C8").Value = 0.3010859938926274
$ws.Range("D8").Value = 0.0007772995109505132
$ws.Range("E8").Value = 0.08751226772843168
$ws.Range("F8").Value = 6.484190562750712
$ws.Range("G8").Value = 0.002742928345600019
$ws.Range("J8").Value = 0.287268135531356
$ws.Range("K8").Value = 5.797244379032861
$ws.Range("L8").Value = 0.05400440540528084
$ws.Range("M8").Value = 1.314777624259065
$ws.Range("N8").Value = 3.634478224043875

$ws.Range("C9").Value = 0.3044979796193843
$ws.Range("D9").Value = 0.0009968044925656727
$ws.Range("E9").Value = 0.08696127162288114
$ws.Range("F9").Value = 6.527135071038543
$ws.Range("G9").Value = 0.00272729203788713
$ws.Range("J9").Value = 0.2851019790195437
$ws.Range("K9").Value = 6.098560286447992
$ws.Range("L9").Value = 0.05419851507762807
$ws.Range("M9").Value = 1.362127334269374
$ws.Range("N9").Value = 3.581231655616335

$ws.Range("C10").Value = 0.3076568594369746
$ws.Range("D10").Value = 0.001163537097077949
$ws.Range("E10").Value = 0.08666956459133957
$ws.Range("F10").Value = 6.573004474524907
$ws.Range("G10").Value = 0.002716830062840956
$ws.Range("J10").Value = 0.2839539197170211
$ws.Range("K10").Value = 6.336655626736786
$ws.Range("L10").Value = 0.05434001879397865
$ws.Range("M10").Value = 1.400410167868046
$ws.Range("N10").Value = 3.546592906435407

$ws.Range("C11").Value = 0.3092362052614845
$ws.Range("D11").Value = 0.001240903852441022
$ws.Range("E11").Value = 0.0865612945734533
$ws.Range("F11").Value = 6.59701261134191
$ws.Range("G11").Value = 0.002712290712136978
$ws.Range("J11").Value = 0.2835280539133578
$ws.Range("K11").Value = 6.44865521876369
$ws.Range("L11").Value = 0.05440421602725287
$ws.Range("M11").Value = 1.41859308701676
$ws.Range("N11").Value = 3.531813928812142

$ws.Range("C12").Value = 0.3098547800186395
$ws.Range("D12").Value = 0.001270446490158506
$ws.Range("E12").Value = 0.08652379795106668
$ws.Range("F12").Value = 6.606558141225321
$ws.Range("G12").Value = 0.002710603183655292
$ws.Range("J12").Value = 0.2833806616287404
$ws.Range("K12").Value = 6.491600847204211
$ws.Range("L12").Value = 0.0544285059869658
$ws.Range("M12").Value = 1.425589479868023
$ws.Range("N12").Value = 3.526358750242792

$ws.Range("C13").Value = 0.3097206459879231
$ws.Range("D13").Value = 0.00126407259828909
$ws.Range("E13").Value = 0.08653171788426661
$ws.Range("F13").Value = 6.6044821029231
$ws.Range("G13").Value = 0.002710965228572986
$ws.Range("J13").Value = 0.2834117879346607
$ws.Range("K13").Value = 6.482327941467474
$ws.Range("L13").Value = 0.05442327552718051
$ws.Range("M13").Value = 1.424077740330063
$ws.Range("N13").Value = 3.527527326319344

$ws.Range("C14").Value = 0.3092866844733777
$ws.Range("D14").Value = 0.001243329265433246
$ws.Range("E14").Value = 0.08655813955700964
$ws.Range("F14").Value = 6.597788811834761
$ws.Range("G14").Value = 0.002712151249453396
$ws.Range("J14").Value = 0.2835156497311431
$ws.Range("K14").Value = 6.452177663383168
$ws.Range("L14").Value = 0.05440621475237783
$ws.Range("M14").Value = 1.419166458625227
$ws.Range("N14").Value = 3.53136229258574

$ws.Range("C15").Value = 0.3090235428604586
$ws.Range("D15").Value = 0.001230656161006038
$ws.Range("E15").Value = 0.08657477948157144
$ws.Range("F15").Value = 6.593748195290374
$ws.Range("G15").Value = 0.002712881808163391
$ws.Range("J15").Value = 0.283581075265019
$ws.Range("K15").Value = 6.433779367490217
$ws.Range("L15").Value = 0.05439576205648411
$ws.Range("M15").Value = 1.416172616412865
$ws.Range("N15").Value = 3.533729739905326

$ws.Range("C16").Value = 0.3075565134088549
$ws.Range("D16").Value = 0.001158513871864386
$ws.Range("E16").Value = 0.08667713100697583
$ws.Range("F16").Value = 6.57149891056406
$ws.Range("G16").Value = 0.002717131128120818
$ws.Range("J16").Value = 0.2839836917709135
$ws.Range("K16").Value = 6.329410703800818
$ws.Range("L16").Value = 0.05433582027165684
$ws.Range("M16").Value = 1.399237364876186
$ws.Range("N16").Value = 3.547578486569506

$ws.Range("C17").Value = 0.3066930268295067
$ws.Range("D17").Value = 0.001114665135439452
$ws.Range("E17").Value = 0.0867461700254708
$ws.Range("F17").Value = 6.558656079307156
$ws.Range("G17").Value = 0.0027197941230288
$ws.Range("J17").Value = 0.2842553810163508
$ws.Range("K17").Value = 6.266331066282021
$ws.Range("L17").Value = 0.05429900688516742
$ws.Range("M17").Value = 1.389045153552956
$ws.Range("N17").Value = 3.556325257696059

$ws.Range("C18").Value = 0.3062097688054166
$ws.Range("D18").Value = 0.001089585855396535
$ws.Range("E18").Value = 0.08678817939509287
$ws.Range("F18").Value = 6.551564809334536
$ws.Range("G18").Value = 0.00272134651229935
$ws.Range("J18").Value = 0.2844207215346515
$ws.Range("K18").Value = 6.230396329349901
$ws.Range("L18").Value = 0.05427781617353045
$ws.Range("M18").Value = 1.383255107427473
$ws.Range("N18").Value = 3.561448219912464

$ws.Range("C19").Value = 0.3060484456405845
$ws.Range("D19").Value = 0.001081117902408835
$ws.Range("E19").Value = 0.08680279837263249
$ws.Range("F19").Value = 6.549214523331841
$ws.Range("G19").Value = 0.002721875686477635
$ws.Range("J19").Value = 0.2844782607546961
$ws.Range("K19").Value = 6.218288942026732
$ws.Range("L19").Value = 0.05427063836219581
$ws.Range("M19").Value = 1.381307090458122
$ws.Range("N19").Value = 3.563198558244295

$ws.Range("C20").Value = 0.3067835596678776
$ws.Range("D20").Value = 0.001119318084350951
$ws.Range("E20").Value = 0.08673858274078761
$ws.Range("F20").Value = 6.559992611561285
$ws.Range("G20").Value = 0.002719508501150273
$ws.Range("J20").Value = 0.2842255202213551
$ws.Range("K20").Value = 6.273010063952881
$ws.Range("L20").Value = 0.05430292741862175
$ws.Range("M20").Value = 1.390122651299876
$ws.Range("N20").Value = 3.555384617681

$ws.Range("C21").Value = 0.3094135925000217
$ws.Range("D21").Value = 0.00124941520394195
$ws.Range("E21").Value = 0.08655028388233532
$ws.Range("F21").Value = 6.599742448485699
$ws.Range("G21").Value = 0.002711802034914088
$ws.Range("J21").Value = 0.2834847663636921
$ws.Range("K21").Value = 6.461019016003547
$ws.Range("L21").Value = 0.05441122642352436
$ws.Range("M21").Value = 1.420606006820393
$ws.Range("N21").Value = 3.530232030258375

$ws.Range("C22").Value = 0.3112520439688495
$ws.Range("D22").Value = 0.001335885211876331
$ws.Range("E22").Value = 0.08644763272604905
$ws.Range("F22").Value = 6.628369653425267
$ws.Range("G22").Value = 0.002706948497615658
$ws.Range("J22").Value = 0.2830815151029284
$ws.Range("K22").Value = 6.587006848466899
$ws.Range("L22").Value = 0.054481890990143
$ws.Range("M22").Value = 1.441175390134319
$ws.Range("N22").Value = 3.514617255419054

$ws.Range("C23").Value = 0.3102598729210371
$ws.Range("D23").Value = 0.001289593347774343
$ws.Range("E23").Value = 0.0865005550173521
$ws.Range("F23").Value = 6.612847650099127
$ws.Range("G23").Value = 0.00270952223199572
$ws.Range("J23").Value = 0.2832893333629016
$ws.Range("K23").Value = 6.519478780879808
$ws.Range("L23").Value = 0.05444418490020109
$ws.Range("M23").Value = 1.430137776041462
$ws.Range("N23").Value = 3.522875562477381

$ws.Range("C24").Value = 0.30674258874744
$ws.Range("D24").Value = 0.00111721408303822
$ws.Range("E24").Value = 0.08674200573166146
$ws.Range("F24").Value = 6.559387455230024
$ws.Range("G24").Value = 0.002719637564156997
$ws.Range("J24").Value = 0.2842389918100352
$ws.Range("K24").Value = 6.269989459653914
$ws.Range("L24").Value = 0.05430115502751143
$ws.Range("M24").Value = 1.389635297842474
$ws.Range("N24").Value = 3.555809587270602

$ws.Range("C25").Value = 0.3034606433566438
$ws.Range("D25").Value = 0.0009365730433170683
$ws.Range("E25").Value = 0.08709042304678505
$ws.Range("F25").Value = 6.513013359902033
$ws.Range("G25").Value = 0.002731340978344265
$ws.Range("J25").Value = 0.2856101504405117
$ws.Range("K25").Value = 6.014127468450113
$ws.Range("L25").Value = 0.05414621918893303
$ws.Range("M25").Value = 1.362127334269374
$ws.Range("N25").Value = 3.59485125094146
